# Generate Report for Handoff
# Adds a new row (row 9) for file c75b0bf8-5cfd-4387-a954-9685ab5a0fea.md
# to the Overview sheet and to the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# OLE color value for the existing hyperlink font color (#6495ED, stored BGR for Excel.Font.Color)
$hyperlinkColor = 15570276

$fileId = "c75b0bf8-5cfd-4387-a954-9685ab5a0fea"
$xlfHash = "15140a5e94629eb678e55ae3f0fc69efb766c8b8"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a95928cd9b57ab0032fd78d397931d5663c0eeb/e2e/$fileId.md",
    "",
    "",
    "$fileId.md"
)
$wsOverview.Range("A9").Font.Underline = $true
$wsOverview.Range("A9").Font.Color = $hyperlinkColor

$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-37-09 10:37:50"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a95928cd9b57ab0032fd78d397931d5663c0eeb/e2e/$fileId.md",
    "",
    "",
    "$fileId.md"
)
$wsZhCn.Range("A9").Font.Underline = $true
$wsZhCn.Range("A9").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a95928cd9b57ab0032fd78d397931d5663c0eeb/e2e/$fileId.md",
    "",
    "",
    ".md"
)
$wsZhCn.Range("B9").Font.Underline = $true
$wsZhCn.Range("B9").Font.Color = $hyperlinkColor

$wsZhCn.Range("C9").Value = "Ready for handoff"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c59e84ee591ab800862925d20c9c3c707304264/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$fileId.$xlfHash.zh-cn.xlf",
    "",
    "",
    "$fileId.$xlfHash.zh-cn.xlf"
)
$wsZhCn.Range("D9").Font.Underline = $true
$wsZhCn.Range("D9").Font.Color = $hyperlinkColor

$wsZhCn.Range("E9").Value = "2016-03-09 10:37:41"
$wsZhCn.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I9").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a95928cd9b57ab0032fd78d397931d5663c0eeb/e2e/$fileId.md",
    "",
    "",
    "$fileId.md"
)
$wsDeDe.Range("A9").Font.Underline = $true
$wsDeDe.Range("A9").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8a95928cd9b57ab0032fd78d397931d5663c0eeb/e2e/$fileId.md",
    "",
    "",
    ".md"
)
$wsDeDe.Range("B9").Font.Underline = $true
$wsDeDe.Range("B9").Font.Color = $hyperlinkColor

$wsDeDe.Range("C9").Value = "Ready for handoff"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eff199f8ae6ac70cba653a5b132fbc873b9afb19/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$fileId.$xlfHash.de-de.xlf",
    "",
    "",
    "$fileId.$xlfHash.de-de.xlf"
)
$wsDeDe.Range("D9").Font.Underline = $true
$wsDeDe.Range("D9").Font.Color = $hyperlinkColor

$wsDeDe.Range("E9").Value = "2016-03-09 10:37:50"
$wsDeDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I9").Value = "Include"

Write-Host "Row 9 added to Overview, zh-cn, de-de sheets for $fileId"
